$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 15) to the "log" sheet, matching the layout of
# the existing rows. Columns B and F:N hold values that look like plain
# integers (meter codes, readings, year/month/day/hour/minute); format them
# as Text first so Excel stores them as strings rather than silently
# re-typing them as numbers, then restore the default "Normal" cell style
# so no stray per-cell number formatting is left on the new row.
$ws.Range("B15").NumberFormat = "@"
$ws.Range("F15:N15").NumberFormat = "@"

$ws.Range("A15").Value = "Kanat"
$ws.Range("B15").Value = "9913"
$ws.Range("C15").Value = "Байсеркешов А"
$ws.Range("D15").Value = "Толе би (Комешбулак)"
$ws.Range("E15").Value = "21/1"
$ws.Range("F15").Value = "10345"
$ws.Range("G15").Value = "1860"
$ws.Range("H15").Value = "1866"
$ws.Range("I15").Value = "6"
$ws.Range("J15").Value = "2025"
$ws.Range("K15").Value = "5"
$ws.Range("L15").Value = "22"
$ws.Range("M15").Value = "20"
$ws.Range("N15").Value = "34"

$ws.Range("B15").Style = "Normal"
$ws.Range("F15:N15").Style = "Normal"
